# Edited the statement of the problemo
#
# 1. Move the "_GoBack" bookmark from its old spot (inside the Objectives
#    section, between "...verification of" and " the servers' authenticity...")
#    to the very start of the "Statement of the Problem" paragraph.
# 2. Rewrite the "Statement of the Problem" paragraph so the question is a
#    single run with the new wording (adds "eliminate the possibility of
#    secret servers to" before "secure the transmission...").

$d = $word.ActiveDocument

# --- Locate the "Statement of the Problem" question paragraph -------------
$targetPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "How can the Philippine*") {
        $targetPara = $p
        break
    }
}

$parStart = $targetPara.Range.Start
$parEnd = $targetPara.Range.End

# --- Step 1: remove the old "_GoBack" bookmark -----------------------------
$oldBookmark = $d.Bookmarks("_GoBack")
$oldBookmark.Delete()

# --- Step 2: re-create it, collapsed, at the start of the question --------
$bookmarkRange = $d.Range($parStart, $parStart)
$d.Bookmarks.Add("_GoBack", $bookmarkRange)

# --- Step 3: rewrite the question text as a single run --------------------
# (exclude the trailing paragraph mark from the replaced range)
$questionRange = $d.Range($parStart, $parEnd - 1)
$questionRange.Text = "How can the Philippine automated election system eliminate the possibility of secret servers to secure the transmission of election returns on the server-level?"
